$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing data row (row 2) values before we overwrite row 1 with them.
$make   = $ws.Range("B2").Value2
$volume = $ws.Range("C2").Value2
$owner  = $ws.Range("D2").Value2
$regDate   = $ws.Range("E2").Value2
$regReason = $ws.Range("F2").Value2
$price     = $ws.Range("G2").Value2

# Move the data from row 2 up into row 1, replacing the old header labels.
$ws.Range("B1").Value = $make
$ws.Range("C1").Value = $volume
$ws.Range("D1").Value = $owner
$ws.Range("E1").Value = $regDate
$ws.Range("F1").Value = $regReason
$ws.Range("G1").Value = $price

# Remove the now-duplicated data row (and the old A2 serial id) entirely.
$ws.Rows(2).Delete()
